# Add a new "Producto" column (D) to the labels/dispensing export sheet,
# holding the product name for every data row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$productName = "jabon liquido - mawie (500 ml)"
$lastRow = 6

# New header
$ws.Cells.Item(1, 4).Value = "Producto"

# New value for each existing data row
for ($row = 2; $row -le $lastRow; $row++) {
    $ws.Cells.Item($row, 4).Value = $productName
}
